$wb = $excel.ActiveWorkbook

# --- Summary sheet: update Total Trades and Win Rate % ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 64
$summary.Range("B9").Value = 32.81

# --- Strategy Status sheet: update MarketMaking row (Trades, Win Rate %) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 64
$status.Range("G4").Value = 32.81

# --- Append new closed trade (#64) to "All Trades" and "MarketMaking" sheets ---
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 65

    $ws.Cells.Item($row, 1).Value = 64

    # Force column B to stay text so the date-like string "2026-02-17" is not
    # auto-converted into a date serial number by Excel.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "15:44:37"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.97
    $ws.Cells.Item($row, 7).Value = 0.97
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.38
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
